# Updated code for the variation 5 and 6.
#
# - All locale sheets except pt_pt: the template-variable code in C1
#   changes from "index-var2" to "index-var5".
# - pt_pt (variation 6): C1 changes from "index-var2" to "index-5".
# - de_de gains a new row (13) duplicating row 2 (the "security" entry).
# - ru_ru becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

$varSheets = @("de_de","ja_jp","zh_cn","ru_ru","ar_ae","fr_fr","es_es","en_au")
foreach ($name in $varSheets) {
  $ws = $wb.Worksheets.Item($name)
  $ws.Range("C1").Value2 = "index-var5"
}

$ptws = $wb.Worksheets.Item("pt_pt")
$ptws.Range("C1").Value2 = "index-5"

# de_de: append a new row 13 that duplicates row 2.
$de = $wb.Worksheets.Item("de_de")
$a2 = $de.Range("A2").Value2
$b2 = $de.Range("B2").Value2
$c2 = $de.Range("C2").Value2
$d2 = $de.Range("D2").Value2
$de.Range("A13").Value2 = $a2
$de.Range("B13").Value2 = $b2
$de.Range("C13").Value2 = $c2
$de.Range("D13").Value2 = $d2

# Restore the per-sheet cursor position (last selected cell) recorded in
# each sheet's view state when the workbook was last saved.
$selections = @{
  "de_de" = "A16";
  "ja_jp" = "B17";
  "zh_cn" = "B14";
  "ru_ru" = "A20";
  "pt_pt" = "B16";
  "ar_ae" = "B16";
  "fr_fr" = "B15";
  "es_es" = "B16";
  "en_au" = "A18";
}
foreach ($name in $selections.Keys) {
  $sheet = $wb.Worksheets.Item($name)
  $sheet.Range($selections[$name]).Select() | Out-Null
}

# ru_ru is the sheet that was active when the workbook was last saved.
$ru = $wb.Worksheets.Item("ru_ru")
$ru.Activate() | Out-Null
